{"js": "// Update the worksheet date and every two-digit multiplication problem's\n// text with the newly generated values, preserving all run/paragraph\n// formatting (fonts, sizes, alignment, etc).\n\n// 1) Update the date heading paragraph.\nconst body = context.document.body;\nconst dateResults = body.search(\"2023-10-12 Thursday\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2023-10-13 Friday\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Update every multiplication-fact cell in the first table. The table is\n// a 20-row x 5-column grid where only every 5th row (0, 4, 9, 14, 19) holds\n// visible equations; the rest are blank spacer rows. Map old -> new text by\n// exact (unique) match so formatting-bearing cells are rewritten in place\n// via Table.values (per-cell Range.Text set), not by touching XML directly.\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"values\");\nawait context.sync();\n\nconst replacements = new Map([\n  [\"76\u00d781=6156\", \"11\u00d712=132\"],\n  [\"98\u00d716=1568\", \"30\u00d726=780\"],\n  [\"35\u00d743=1505\", \"41\u00d735=1435\"],\n  [\"95\u00d715=1425\", \"62\u00d756=3472\"],\n  [\"96\u00d739=3744\", \"20\u00d746=920\"],\n  [\"54\u00d778=4212\", \"35\u00d734=1190\"],\n  [\"92\u00d769=6348\", \"63\u00d718=1134\"],\n  [\"36\u00d727=972\", \"69\u00d726=1794\"],\n  [\"14\u00d724=336\", \"57\u00d767=3819\"],\n  [\"23\u00d774=1702\", \"37\u00d769=2553\"],\n  [\"94\u00d722=2068\", \"68\u00d730=2040\"],\n  [\"98\u00d779=7742\", \"72\u00d768=4896\"],\n  [\"50\u00d776=3800\", \"47\u00d796=4512\"],\n  [\"92\u00d721=1932\", \"33\u00d721=693\"],\n  [\"13\u00d793=1209\", \"42\u00d794=3948\"],\n  [\"50\u00d750=2500\", \"57\u00d764=3648\"],\n  [\"78\u00d747=3666\", \"95\u00d778=7410\"],\n  [\"81\u00d732=2592\", \"65\u00d754=3510\"],\n  [\"91\u00d774=6734\", \"14\u00d715=210\"],\n  [\"82\u00d782=6724\", \"88\u00d795=8360\"],\n  [\"86\u00d768=5848\", \"42\u00d753=2226\"],\n  [\"59\u00d773=4307\", \"98\u00d767=6566\"],\n  [\"40\u00d728=1120\", \"74\u00d784=6216\"],\n  [\"96\u00d729=2784\", \"34\u00d763=2142\"],\n  [\"38\u00d777=2926\", \"32\u00d750=1600\"],\n]);\n\nif (!table.isNullObject) {\n  for (let r = 0; r < table.values.length; r++) {\n    const row = table.values[r];\n    for (let c = 0; c < row.length; c++) {\n      const oldText = row[c];\n      if (replacements.has(oldText)) {\n        table.getCell(r, c).value = replacements.get(oldText);\n      }\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every two-digit multiplication problem's\n# text with the newly generated values, preserving all run/paragraph\n# formatting (fonts, sizes, alignment, etc) by using Find/Replace (which\n# only rewrites the matched text, not the surrounding run properties).\n\n$d = $word.ActiveDocument\n\n# Map of old exact text -> new exact text. All old values are unique across\n# the whole document, so a plain Replace-All of the whole match is safe and\n# unambiguous (no partial/substring collisions between entries).\n$replacements = [ordered]@{\n    \"2023-10-12 Thursday\" = \"2023-10-13 Friday\"\n    \"76\u00d781=6156\" = \"11\u00d712=132\"\n    \"98\u00d716=1568\" = \"30\u00d726=780\"\n    \"35\u00d743=1505\" = \"41\u00d735=1435\"\n    \"95\u00d715=1425\" = \"62\u00d756=3472\"\n    \"96\u00d739=3744\" = \"20\u00d746=920\"\n    \"54\u00d778=4212\" = \"35\u00d734=1190\"\n    \"92\u00d769=6348\" = \"63\u00d718=1134\"\n    \"36\u00d727=972\" = \"69\u00d726=1794\"\n    \"14\u00d724=336\" = \"57\u00d767=3819\"\n    \"23\u00d774=1702\" = \"37\u00d769=2553\"\n    \"94\u00d722=2068\" = \"68\u00d730=2040\"\n    \"98\u00d779=7742\" = \"72\u00d768=4896\"\n    \"50\u00d776=3800\" = \"47\u00d796=4512\"\n    \"92\u00d721=1932\" = \"33\u00d721=693\"\n    \"13\u00d793=1209\" = \"42\u00d794=3948\"\n    \"50\u00d750=2500\" = \"57\u00d764=3648\"\n    \"78\u00d747=3666\" = \"95\u00d778=7410\"\n    \"81\u00d732=2592\" = \"65\u00d754=3510\"\n    \"91\u00d774=6734\" = \"14\u00d715=210\"\n    \"82\u00d782=6724\" = \"88\u00d795=8360\"\n    \"86\u00d768=5848\" = \"42\u00d753=2226\"\n    \"59\u00d773=4307\" = \"98\u00d767=6566\"\n    \"40\u00d728=1120\" = \"74\u00d784=6216\"\n    \"96\u00d729=2784\" = \"34\u00d763=2142\"\n    \"38\u00d777=2926\" = \"32\u00d750=1600\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
